$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new feed-log rows (24 and 25) after the existing data (row 23)
$newRows = @(
    @(23, 1, "2024-06-15 09:12:29", 200, 2),
    @(24, 2, "2024-06-15 09:12:30", 200, 0)
)

$r = 24
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
